$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 4526.3887
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 5127.4194
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 15382.2582
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -17598.2582
# Row 132
$ws.Range("H132").Value = 29457588
$ws.Range("I132").Value = 32804840
$ws.Range("J132").Value = 1760
$ws.Range("K132").Value = 98414520
$ws.Range("L132").Value = 5280
$ws.Range("M132").Value = -98411990
$ws.Range("N132").Value = -10340
# Row 135
$ws.Range("H135").Value = 399717.38
$ws.Range("I135").Value = 3052.1667
$ws.Range("K135").Value = 27469.5003
$ws.Range("M135").Value = -24934.5003
# Row 137
$ws.Range("H137").Value = 332272.28
$ws.Range("I137").Value = 542944.4399999999
$ws.Range("J137").Value = 1216.0714
$ws.Range("K137").Value = 1628833.32
$ws.Range("L137").Value = 3648.2142
$ws.Range("M137").Value = -1626283.32
$ws.Range("N137").Value = -8748.2142
# Row 138
$ws.Range("H138").Value = 1907.284
$ws.Range("I138").Value = 1271.5238
$ws.Range("J138").Value = 2487.761
$ws.Range("K138").Value = 3814.5714
$ws.Range("L138").Value = 7463.282999999999
$ws.Range("M138").Value = 1325.4286
$ws.Range("N138").Value = -17743.283
# Row 141
$ws.Range("H141").Value = 2743.561
$ws.Range("I141").Value = 2180.3381
$ws.Range("J141").Value = 5479.2144
$ws.Range("K141").Value = 6541.0143
$ws.Range("L141").Value = 16437.6432
$ws.Range("M141").Value = -1361.0143
$ws.Range("N141").Value = -26797.6432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1859.34
$ws.Range("I32").Value = 1732.4731
$ws.Range("J32").Value = 3544.8572
$ws.Range("K32").Value = 1732.4731
$ws.Range("L32").Value = 3544.8572
$ws.Range("M32").Value = -1445.4731
$ws.Range("N32").Value = -4118.8572
# Row 74
$ws.Range("H74").Value = 703.3333
$ws.Range("I74").Value = 507.45947
$ws.Range("J74").Value = 1018.43475
$ws.Range("K74").Value = 507.45947
$ws.Range("L74").Value = 1018.43475
$ws.Range("M74").Value = 366.54053
$ws.Range("N74").Value = -2766.43475
# Row 77
$ws.Range("H77").Value = 703.3333
$ws.Range("I77").Value = 507.45947
$ws.Range("J77").Value = 1018.43475
$ws.Range("K77").Value = 2537.29735
$ws.Range("L77").Value = 5092.17375
$ws.Range("M77").Value = 1830.70265
$ws.Range("N77").Value = -13828.17375
# Row 132
$ws.Range("H132").Value = 2718702.8
$ws.Range("I132").Value = 3572427.5
$ws.Range("J132").Value = 2306.182
$ws.Range("K132").Value = 10717282.5
$ws.Range("L132").Value = 6918.545999999999
$ws.Range("M132").Value = -10714752.5
$ws.Range("N132").Value = -11978.546

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 50001450
$ws.Range("I107").Value = 62501468
$ws.Range("K107").Value = 62501468
$ws.Range("M107").Value = -62499548
# Row 134
$ws.Range("H134").Value = 8561104
$ws.Range("I134").Value = 10117378
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 30352134
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -30349599
$ws.Range("N134").Value = -9870

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8330.5
$ws.Range("I31").Value = 1133.1143
$ws.Range("J31").Value = 44317.43
$ws.Range("K31").Value = 1133.1143
$ws.Range("L31").Value = 44317.43
$ws.Range("M31").Value = -838.1143
$ws.Range("N31").Value = -44907.43
# Row 34
$ws.Range("H34").Value = 8330.5
$ws.Range("I34").Value = 1133.1143
$ws.Range("J34").Value = 44317.43
$ws.Range("K34").Value = 1133.1143
$ws.Range("L34").Value = 44317.43
$ws.Range("M34").Value = -931.1143
$ws.Range("N34").Value = -44721.43
# Row 58
$ws.Range("H58").Value = 3428598.5
$ws.Range("I58").Value = 4961467.5
$ws.Range("J58").Value = 9122.154
$ws.Range("K58").Value = 4961467.5
$ws.Range("L58").Value = 9122.154
$ws.Range("M58").Value = -4961264.5
$ws.Range("N58").Value = -9528.154
# Row 134
$ws.Range("H134").Value = 29167924
$ws.Range("I134").Value = 44643972
$ws.Range("J134").Value = 3677958.2
$ws.Range("K134").Value = 133931916
$ws.Range("L134").Value = 11033874.6
$ws.Range("M134").Value = -133929381
$ws.Range("N134").Value = -11038944.6
# Row 136
$ws.Range("H136").Value = 3428598.5
$ws.Range("I136").Value = 4961467.5
$ws.Range("J136").Value = 9122.154
$ws.Range("K136").Value = 14884402.5
$ws.Range("L136").Value = 27366.462
$ws.Range("M136").Value = -14881852.5
$ws.Range("N136").Value = -32466.462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 10416.55
$ws.Range("I4").Value = 152.5
$ws.Range("J4").Value = 25812.625
$ws.Range("K4").Value = 457.5
$ws.Range("L4").Value = 77437.875
$ws.Range("M4").Value = -345.5
$ws.Range("N4").Value = -77661.875
# Row 17
$ws.Range("H17").Value = 142993
$ws.Range("I17").Value = 500025.5
$ws.Range("J17").Value = 180
$ws.Range("K17").Value = 1500076.5
$ws.Range("L17").Value = 540
$ws.Range("M17").Value = -1499907.5
$ws.Range("N17").Value = -878

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 19610910
$ws.Range("I132").Value = 32259606
$ws.Range("J132").Value = 5432.3
$ws.Range("K132").Value = 96778818
$ws.Range("L132").Value = 16296.9
$ws.Range("M132").Value = -96776288
$ws.Range("N132").Value = -21356.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4701.125
$ws.Range("I16").Value = 550.8095
$ws.Range("J16").Value = 33753.332
$ws.Range("K16").Value = 550.8095
$ws.Range("L16").Value = 33753.332
$ws.Range("M16").Value = -380.8095
$ws.Range("N16").Value = -34093.332
# Row 132
$ws.Range("H132").Value = 4446356.5
$ws.Range("I132").Value = 5715561.5
$ws.Range("J132").Value = 4140.3
$ws.Range("K132").Value = 17146684.5
$ws.Range("L132").Value = 12420.9
$ws.Range("M132").Value = -17144154.5
$ws.Range("N132").Value = -17480.9
# Row 136
$ws.Range("H136").Value = 3501.3333
$ws.Range("I136").Value = 3734.152
$ws.Range("J136").Value = 2527.7273
$ws.Range("K136").Value = 11202.456
$ws.Range("L136").Value = 7583.1819
$ws.Range("M136").Value = -8652.456
$ws.Range("N136").Value = -12683.1819

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
# Row 132
$ws.Range("H132").Value = 8501477
$ws.Range("I132").Value = 5129363.5
$ws.Range("J132").Value = 11423975
$ws.Range("K132").Value = 15388090.5
$ws.Range("L132").Value = 34271925
$ws.Range("M132").Value = -15385560.5
$ws.Range("N132").Value = -34276985
# Row 136
$ws.Range("H136").Value = 12605081
$ws.Range("I136").Value = 8710216
$ws.Range("J136").Value = 18519506
$ws.Range("K136").Value = 26130648
$ws.Range("L136").Value = 55558518
$ws.Range("M136").Value = -26128098
$ws.Range("N136").Value = -55563618
